# Added check for additional samples and expanded model.
#
# The "attributes" sheet (sheet3.xml) describes the data model's attribute
# list. A new attribute row - "sample_id2" / "Combined virtual samples and
# novelomics IDs" - is inserted right after the existing "sample_id" row
# (entity = solverdportal_experiments), pushing every following row down by
# one. The "attributes" tab also becomes the active tab/sheet of the
# workbook (it was previously "solverdportal_experiment_errors").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# Insert a new row above the current row 4 (participant_id), shifting
# everything below it down by one row.
$ws.Rows("4").Insert()

# Populate the newly inserted row 4 with the new attribute definition.
$ws.Cells.Item(4, 1).Value = "solverdportal_experiments"
$ws.Cells.Item(4, 2).Value = "sample_id2"
$ws.Cells.Item(4, 3).Value = "Combined virtual samples and novelomics IDs"
$ws.Cells.Item(4, 4).Value = "string"

# Make "attributes" the active sheet/tab and select E4, matching the
# worked-in selection left behind after adding the row.
$ws.Activate() | Out-Null
$ws.Range("E4").Select() | Out-Null
